$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2-13 get cyclically re-shuffled: the new content of each row
# equals the old content of another row (a full 12-cycle permutation).
# Mapping: new row -> source row (both 2..13), derived from the diff.
$mapping = @{
    2  = 3
    3  = 4
    4  = 5
    5  = 6
    6  = 11
    7  = 13
    8  = 2
    9  = 7
    10 = 8
    11 = 9
    12 = 10
    13 = 12
}

# Columns used anywhere in rows 2-13 (A..AY span, but only populated ones matter)
$cols = @("A","B","C","D","E","F","G","H","I","J","K","N","P","Q","R","S","T","U","V","W","Y","Z","AA","AB","AD","AE","AF","AG","AT","AW","AX","AY")

# Snapshot the current (before) values of every relevant cell for rows 2-13
$snapshot = @{}
foreach ($row in 2..13) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $addr = "$col$row"
        $rowVals[$col] = $ws.Range($addr).Value()
    }
    $snapshot[$row] = $rowVals
}

# Now write back: each destination row gets the snapshot of its source row
foreach ($row in 2..13) {
    $src = $mapping[$row]
    $rowVals = $snapshot[$src]
    foreach ($col in $cols) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $rowVals[$col]
    }
}
